# Apply the 29 Jan 2024 cryptos list refresh (prices + 1h volume deltas),
# including the PancakeSwap/ImmutableX and Hedera/Celestia rank swaps.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.270.96"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.269.33"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("E7").Value = "  -0.80%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.42%  "
$ws.Range("D14").Value = "2.623.97"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "2.295.25"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.792"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "42.128.35"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.95%  "
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.01"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.79%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.92%  "
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.38%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0738"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  -3.92%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.77%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.25%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.942.90"
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.20"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.59%  "
